{"js": "// The document is a date header paragraph followed by a 5-column table\n// whose populated rows contain \"A\u00f7B=C, D\" style division answers.\n// The edit replaces the header date/weekday and every populated answer\n// cell's text with a new value, in document order (paragraphs, including\n// those inside table cells, are visited top-to-bottom, left-to-right).\nconst replacements = [\n  \"2024-02-02 Friday\",\n  \"696\u00f72=348, 0\",\n  \"154\u00f74=38, 2\",\n  \"423\u00f74=105, 3\",\n  \"212\u00f75=42, 2\",\n  \"835\u00f74=208, 3\",\n  \"699\u00f75=139, 4\",\n  \"293\u00f75=58, 3\",\n  \"497\u00f72=248, 1\",\n  \"269\u00f77=38, 3\",\n  \"354\u00f75=70, 4\",\n  \"331\u00f78=41, 3\",\n  \"539\u00f73=179, 2\",\n  \"187\u00f73=62, 1\",\n  \"260\u00f75=52, 0\",\n  \"780\u00f73=260, 0\",\n  \"321\u00f72=160, 1\",\n  \"605\u00f74=151, 1\",\n  \"114\u00f79=12, 6\",\n  \"586\u00f76=97, 4\",\n  \"498\u00f72=249, 0\",\n  \"240\u00f73=80, 0\",\n  \"296\u00f79=32, 8\",\n  \"418\u00f73=139, 1\",\n  \"978\u00f77=139, 5\",\n  \"646\u00f77=92, 2\",\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet idx = 0;\nfor (const paragraph of paragraphs.items) {\n  if (idx >= replacements.length) break;\n  // Skip the blank spacer paragraphs that live in the table's empty rows.\n  if (paragraph.text === \"\") continue;\n  paragraph.insertText(replacements[idx], \"Replace\");\n  idx++;\n}\n\nawait context.sync();\n", "ps1": "# The document is a date header paragraph followed by a 5-column table\n# whose populated rows contain \"A\u00f7B=C, D\" style division answers.\n# The edit replaces the header date/weekday and every populated answer\n# cell's text with a new value, in document order (Document.Paragraphs\n# walks the header paragraph, then each table cell's paragraph, top to\n# bottom / left to right, including the empty spacer rows).\n$d = $word.ActiveDocument\n\n$texts = @(\n  \"2024-02-02 Friday\",\n  \"696\u00f72=348, 0\",\n  \"154\u00f74=38, 2\",\n  \"423\u00f74=105, 3\",\n  \"212\u00f75=42, 2\",\n  \"835\u00f74=208, 3\",\n  \"699\u00f75=139, 4\",\n  \"293\u00f75=58, 3\",\n  \"497\u00f72=248, 1\",\n  \"269\u00f77=38, 3\",\n  \"354\u00f75=70, 4\",\n  \"331\u00f78=41, 3\",\n  \"539\u00f73=179, 2\",\n  \"187\u00f73=62, 1\",\n  \"260\u00f75=52, 0\",\n  \"780\u00f73=260, 0\",\n  \"321\u00f72=160, 1\",\n  \"605\u00f74=151, 1\",\n  \"114\u00f79=12, 6\",\n  \"586\u00f76=97, 4\",\n  \"498\u00f72=249, 0\",\n  \"240\u00f73=80, 0\",\n  \"296\u00f79=32, 8\",\n  \"418\u00f73=139, 1\",\n  \"978\u00f77=139, 5\",\n  \"646\u00f77=92, 2\"\n)\n\n$idx = 0\nforeach ($p in $d.Paragraphs) {\n    if ($idx -ge $texts.Length) { break }\n    $r = $p.Range\n    # Strip the trailing paragraph mark / cell mark before checking for\n    # emptiness so the blank spacer rows between data rows are skipped.\n    $plain = $r.Text.TrimEnd([char]13, [char]7)\n    if ($plain.Length -eq 0) { continue }\n    $r.Text = $texts[$idx]\n    $idx++\n}\n"}
